$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Alchemists bug - update description and change reviser to Jan
$ws.Range("A9").Value = "Alchemists Graphic Pivots are off -> looks weird in champion select --> UI Problem"
$ws.Range("C9").Value = "Jan"

# Row 10: Knight Skill 4 Spear too high - fix text (remove trailing "+"), mark as Fixed, add reviser, add diagonal strike style like rows 2-8
$ws.Range("A10").Value = "Knight Skill 4 Spear too high"
$ws.Range("B10").Value = "Fixed"
$ws.Range("C10").Value = "Sandro"
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Row 11: Knight attack speed too fast - mark as Fixed, add reviser, apply same style
$ws.Range("B11").Value = "Fixed"
$ws.Range("C11").Value = "Sandro"
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Row 13: Fall through not working with controller - add reviser Kyle
$ws.Range("C13").Value = "Kyle"

# Row 14: Fall through causes stick in ground controller - add reviser Kyle
$ws.Range("C14").Value = "Kyle"

# Update selection to A16 (as in target)
$ws.Range("A16").Select()

$wb.Save()
